$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.151.49"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.56"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  +0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.66"
$ws.Range("E5").Value = "  +2.16%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.64"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.326"
$ws.Range("E9").Value = "  +7.53%  "

$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1000"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.081.25"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.828.96"
$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.11"
$ws.Range("E14").Value = "  -3.27%  "

$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.65"
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.117.83"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.69"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.71"
$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.94"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.40"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.85"
$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.49"
$ws.Range("E27").Value = "  -2.08%  "

$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("E29").Value = "  +20.84%  "

$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.333.29"
$ws.Range("E31").Value = "  +37.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.07"
$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0554"
$ws.Range("E33").Value = "  +4.01%  "

$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E36").Value = "  +6.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "93.09"
$ws.Range("E37").Value = "  +1.89%  "

$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0194"
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("E40").Value = "  +2.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.309.92"
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.72"
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.48"
$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.31"
$ws.Range("E45").Value = "  -5.15%  "

$ws.Range("E46").Value = "  -2.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.33"
$ws.Range("E47").Value = "  +5.02%  "

$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.996.63"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0646"
$ws.Range("E51").Value = "  +4.82%  "
